$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.6873773333333334
$ws.Range("H2").Value = 2.062132
$ws.Range("I2").Value = 0.02660947569874856
$ws.Range("J2").Value = 0.02660947569874856
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 36.36934344025023
$ws.Range("R2").Value = 327.324090962252
$ws.Range("S2").Value = 0.01107330826389577
$ws.Range("T2").Value = 0.01107330826389577

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.6873773333333334
$ws.Range("H3").Value = 2.062132
$ws.Range("I3").Value = 0.02660947569874856
$ws.Range("J3").Value = 0.02660947569874856
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 32.48549287034444
$ws.Range("R3").Value = 292.3694358331
$ws.Range("S3").Value = 0.009890799300484627
$ws.Range("T3").Value = 0.009890799300484627

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.6873773333333334
$ws.Range("H4").Value = 2.062132
$ws.Range("I4").Value = 0.02660947569874856
$ws.Range("J4").Value = 0.02660947569874856
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 18.54173365650045
$ws.Range("R4").Value = 166.875602908504
$ws.Range("S4").Value = 0.00564536813436817
$ws.Range("T4").Value = 0.00564536813436817

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.913984666666667
$ws.Range("H5").Value = 17.741954
$ws.Range("I5").Value = 0.2289398029860915
$ws.Range("J5").Value = 0.2289398029860915
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 312.9107245933438
$ws.Range("R5").Value = 2816.196521340094
$ws.Range("S5").Value = 0.09527136276720331
$ws.Range("T5").Value = 0.09527136276720333

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.913984666666667
$ws.Range("H6").Value = 17.741954
$ws.Range("I6").Value = 0.2289398029860915
$ws.Range("J6").Value = 0.2289398029860915
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 279.4952603291055
$ws.Range("R6").Value = 2515.45734296195
$ws.Range("S6").Value = 0.08509741675723494
$ws.Range("T6").Value = 0.08509741675723495

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.913984666666667
$ws.Range("H7").Value = 17.741954
$ws.Range("I7").Value = 0.2289398029860915
$ws.Range("J7").Value = 0.2289398029860915
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 159.5274141586876
$ws.Range("R7").Value = 1435.746727428188
$ws.Range("S7").Value = 0.04857102346165321
$ws.Range("T7").Value = 0.04857102346165322

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.230689
$ws.Range("H8").Value = 57.692067
$ws.Range("I8").Value = 0.7444507213151601
$ws.Range("J8").Value = 0.7444507213151601
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 1017.501594709226
$ws.Range("R8").Value = 9157.514352383037
$ws.Range("S8").Value = 0.3097968715253573
$ws.Range("T8").Value = 0.3097968715253573

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.230689
$ws.Range("H9").Value = 57.692067
$ws.Range("I9").Value = 0.7444507213151601
$ws.Range("J9").Value = 0.7444507213151601
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 908.8434839301916
$ws.Range("R9").Value = 8179.591355371725
$ws.Range("S9").Value = 0.2767139329233591
$ws.Range("T9").Value = 0.2767139329233591

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.230689
$ws.Range("H10").Value = 57.692067
$ws.Range("I10").Value = 0.7444507213151601
$ws.Range("J10").Value = 0.7444507213151601
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 518.7402845244527
$ws.Range("R10").Value = 4668.662560720074
$ws.Range("S10").Value = 0.1579399168664438
$ws.Range("T10").Value = 0.1579399168664438
